# Switch to captest_prototype alternative capacity testing framework;
# include paper and oral presentation.
#
# Inserts a new "Bifaciality" column (H) on the "PVsyst Runs" sheet, plus
# five new trailing columns (sep, dayfirst, date_format, StrucShd,
# BakMismatch), fills in the values for all 9 data rows, and updates the
# view state on sheet1/sheet2 to match the saved selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PVsyst Runs")

# --- Insert the new "Bifaciality" column before the existing "Height" column (H) ---
# (do this first so later column letters - L, M, N, O, P - land in their
# final positions instead of shifting.)
$ws.Columns("H").Insert()

# --- New trailing columns: sep, dayfirst, (Bifaciality header), date_format, StrucShd, BakMismatch ---
$ws.Range("L1").Value = "sep"
$ws.Range("M1").Value = "dayfirst"
$ws.Range("L2").Value = ","
$ws.Range("H1").Value = "Bifaciality"
$ws.Range("N1").Value = "date_format"
$ws.Range("N2").Value = "%m/%d/%y %H:%M"
$ws.Range("O1").Value = "StrucShd"
$ws.Range("P1").Value = "BakMismatch"

# --- Fill "Bifaciality" column values ---
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0.7
$ws.Range("H4").Value = 0.7
$ws.Range("H5").Value = 0.7
$ws.Range("H6").Value = 0.7
$ws.Range("H7").Value = 0
$ws.Range("H8").Value = 0.7
$ws.Range("H9").Value = 0.7

# --- Fill remaining "sep" / "dayfirst" / "date_format" column values ---
$ws.Range("L3:L9").Value = ","
$ws.Range("M2:M9").Value = $false
$ws.Range("N3:N9").Value = "%m/%d/%y %H:%M"

# --- Fill "StrucShd" / "BakMismatch" column values ---
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0

$ws.Range("O3").Value = 0.05
$ws.Range("P3").Value = 0.1

$ws.Range("O4").Value = 0.05
$ws.Range("P4").Value = 0.1

$ws.Range("O5").Value = 0.05
$ws.Range("P5").Value = 0.1

$ws.Range("O6").Value = 0.05
$ws.Range("P6").Value = 0.1

$ws.Range("O7").Value = 0
$ws.Range("P7").Value = 0

$ws.Range("O8").Value = 0.05
$ws.Range("P8").Value = 0.1

$ws.Range("O9").Value = 0.05
$ws.Range("P9").Value = 0.1

# --- Update sheet2 ("Systems") view state: scrolled to D1, selection Y1:Z3 ---
$ws2 = $wb.Worksheets.Item("Systems")
$ws2.Activate()
$ws2.Range("Y1:Z3").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 4

# --- Update sheet1 ("PVsyst Runs") view state: scrolled to B1, active cell O1 ---
$ws.Activate()
$ws.Range("O1").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 2
